$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 232, pushing existing rows 232-238 down to 234-240.
$ws.Rows("232:233").Insert()

# Fill in the new row 232 (copy of the old row 232's static fields, with updated
# date/volume/price data for the new week).
$ws.Range("A232").Value = 2
$ws.Range("B232").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C232").Value = "Coquimbo"
$ws.Range("D232").Value = 45267
$ws.Range("E232").Value = 4
$ws.Range("F232").Value = 100112043
$ws.Range("G232").Value = "Pepino ensalada"
$ws.Range("H232").Value = "Sin especificar"
$ws.Range("I232").Value = "Primera"
$ws.Range("J232").Value = 500
$ws.Range("K232").Value = 14000
$ws.Range("L232").Value = 15000
$ws.Range("M232").Value = 14500
$ws.Range("N232").Value = "$/caja 70 unidades"
$ws.Range("O232").Value = "Provincia de Limarí"
$ws.Range("P232").Value = 207
$ws.Range("Q232").Value = 70
$ws.Range("R232").Value = "Hortaliza"

# Fill in the new row 233 (copy of the old row 233's static fields, with updated
# date/volume/price data for the new week).
$ws.Range("A233").Value = 2
$ws.Range("B233").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C233").Value = "Coquimbo"
$ws.Range("D233").Value = 45267
$ws.Range("E233").Value = 4
$ws.Range("F233").Value = 100112043
$ws.Range("G233").Value = "Pepino ensalada"
$ws.Range("H233").Value = "Sin especificar"
$ws.Range("I233").Value = "Segunda"
$ws.Range("J233").Value = 300
$ws.Range("K233").Value = 11000
$ws.Range("L233").Value = 12000
$ws.Range("M233").Value = 11500
$ws.Range("N233").Value = "$/caja 100 unidades"
$ws.Range("O233").Value = "Provincia de Limarí"
$ws.Range("P233").Value = 115
$ws.Range("Q233").Value = 100
$ws.Range("R233").Value = "Hortaliza"
